$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Edit 1: "ID number*" cell: "2" -> "bug" + "2" (two separate runs) ---
$cell1 = $t.Cell(1, 2)
$pos1 = $cell1.Range.Start
$insertRng1 = $d.Range($pos1, $pos1)
$insertRng1.InsertBefore("bug")
# Toggling a character-formatting property and reverting it keeps the
# newly-inserted "bug" text as its own run instead of merging back into
# the adjacent "2" run (which has identical formatting).
$bugRng = $d.Range($pos1, $pos1 + 3)
$bugRng.Font.Bold = 1
$bugRng.Font.Bold = 0

# --- Edit 2: "Browser (optional)" cell: merge "Google Chrome v" + "ersion 123.0.6312.107" into a single run ---
$cell2 = $t.Cell(8, 2)
$rng2 = $cell2.Range
$start2 = $rng2.Start
$end2 = $rng2.End
# Assigning identical text is a no-op for the engine's run structure, so
# first write a text that differs (adds a trailing sentinel) to force a
# real text replacement, merging the two runs into one, then trim the
# sentinel back off with a second assignment.
$work2 = $d.Range($start2, $end2 - 1)
$work2.Text = "Google Chrome version 123.0.6312.107_"
$work2b = $d.Range($start2, $end2)
$work2b.Text = "Google Chrome version 123.0.6312.107"

# --- Edit 3: "4" + '.press "login"' -> merge into a single run ---
$target3 = $null
foreach ($p in $d.Content.Paragraphs) {
  if ($p.Range.Text -like '*press*') {
    $target3 = $p.Range
  }
}
$start3 = $target3.Start
$end3 = $target3.End
$work3 = $d.Range($start3, $end3 - 1)
$work3.Text = '4.press "login"_'
$work3b = $d.Range($start3, $end3)
$work3b.Text = '4.press "login"'

Write-Host "Cell1:" $t.Cell(1,2).Range.Text
Write-Host "Cell8:" $t.Cell(8,2).Range.Text
Write-Host "Para3:" $d.Range($start3, $end3).Text
